$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Table 1 ("Bill No" details) - simple value updates via Find/Replace. Each
# of these strings is unique in the document so whole-document replace is
# safe.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2020-12-26 13:19:47", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2020-12-28 12:55:40", 2)
$d.Content.Find.Execute("123456", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "123", 2)
$d.Content.Find.Execute("Abcd", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Test", 2)
$d.Content.Find.Execute("07777", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "123", 2)

# ---------------------------------------------------------------------------
# Table 2 ("Rented Item details") - drop the "දින ගණන" and "මුදල" columns,
# widen the remaining four columns, refresh the first item's date, and turn
# the old blank/"Total" row into a second rented-item row.
# ---------------------------------------------------------------------------
$t2 = $d.Tables.Item(2)

# Remove the "දින ගණන" column (3rd column).
$t2.Columns(3).Delete()
# Remove the "මුදල" column (was 6th, now the last/5th column).
$t2.Columns(5).Delete()

# Update the existing item row's rented date.
$d.Content.Find.Execute("2020-12-22", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2020-12-28", 2)

# Replace the old (mostly blank) "Total" row with a new rented-item row.
$t2.Rows(3).Delete()
$newRow = $t2.Rows.Add()
$newRow.Cells(1).Range.Text = "Grinder"
$newRow.Cells(2).Range.Text = "2020-12-28"
$newRow.Cells(3).Range.Text = "1"
$newRow.Cells(4).Range.Text = "4,000.00"

# Widen all remaining columns from 1440 dxa (72pt) to 2160 dxa (108pt).
for ($i = 1; $i -le $t2.Columns.Count; $i++) {
  $t2.Columns($i).Width = 108
}

# ---------------------------------------------------------------------------
# Table 3 ("Payment Details") - keep only the "ගෙවීම්" row, with an updated
# amount; drop the other four rows.
# ---------------------------------------------------------------------------
$t3 = $d.Tables.Item(3)

# Drop the first three rows (හිඟ මුදල, මෙම බිල්පතෙහි වටිනාකම, මුළු මුදල).
$t3.Rows(1).Delete()
$t3.Rows(1).Delete()
$t3.Rows(1).Delete()

# Drop the trailing row (ගෙවිය යුතු වටිනාකම), leaving only ගෙවීම්.
$t3.Rows(2).Delete()

# Update the ගෙවීම් amount.
$t3.Cell(1, 2).Range.Text = "Rs.   2,000.00"
